# Auto update Excel log
# Appends new mmWave PRESENCE_DETECTED log rows (61-67) to the "mmWave" sheet,
# matching the sensor log format already used in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mmWave")

$rows = @(
    @("2026-02-01", "16:02:10", "16:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "16:02:13", "16:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "16:02:23", "16:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "16:02:34", "16:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "16:02:44", "16:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "16:02:55", "16:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "16:03:05", "16:00", "Living Room", "PRESENCE_DETECTED", "Active")
)

$startRow = 61
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    # Column A holds a date-like string ("YYYY-MM-DD"). Excel's COM layer
    # auto-parses that as a real date/time, so force literal text here and
    # then strip the temporary formatting so the cell ends up style-plain,
    # identical to how the existing rows are stored.
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $data[0]
    $cellA.ClearFormats()

    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
}
